# Update the LockedLoans key/value table (column A = index, column B = value)
# to reflect the new loan record, per the authoritative diff.
# NumberFormat is forced to "@" (Text) before assigning any value that
# "looks like" a number, so Excel keeps storing it as a shared string
# (t="s") exactly like the original workbook, instead of silently
# converting it into a numeric cell (which would also strip things like
# trailing/leading zeros). Boolean-looking text ("true"/"false") is
# written with a leading apostrophe so Excel treats it as literal text
# rather than a TRUE/FALSE boolean cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "180000"

$ws.Range("B4").Value = "Detached"
$ws.Range("B5").Value = "DU"
$ws.Range("B6").Value = "Denis"
$ws.Range("B7").Value = "TPO_Y"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "94.31"

$ws.Range("B10").Value = "Blue Island"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "10.0"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "777"

$ws.Range("B13").Value = "'true"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "94.31"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "169750"

$ws.Range("B19").Value = "DEN_20251006_10001_TPO_Y"

$ws.Range("B22").Value = "Null Value"
$ws.Range("B23").Value = "Conventional"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "7.500"

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "1"

$ws.Range("B27").Value = "212_25/30 Yr Fannie Mae Fixed"
$ws.Range("B28").Value = "SingleFamily"

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "175000"

$ws.Range("B32").Value = "12714 Mozart St"

$ws.Range("B34").Value = "'true"

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "169750"

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "169750"

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "60406"

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "212"

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "15"

$ws.Range("B41").Value = "[15]"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "8234.00"

$ws.Range("B49").Value = "T"

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "1.3"
